# #59: fixed data provider RAM values for validation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GSMArena")

# Galaxy S10+ (row 2) actually ships with 8/12GB RAM options, not a flat 12GB RAM
$ws.Range("G2").Value = "8/12GB RAM"

# Galaxy S10 (row 3) only has a single 8GB RAM option, not 8/12GB RAM
$ws.Range("G3").Value = "8GB RAM"
